# Generate Report for Handoff
# Refresh the "Latest Handoff Date(time)" columns for every file that is
# still awaiting handback (i.e. not yet "Handed back: in sync with en-US"
# and not "In Translation"), reflecting a fresh handoff run.

$wb = $excel.ActiveWorkbook

# New handoff timestamps produced by this run, per locale.
$zhcnHandoffTime = "2016-03-24 16:31:17"
$dedeHandoffTime = "2016-03-24 16:31:24"
# Overview's "Latest Handoff Date" reflects the most recent handoff across
# all locales, i.e. the de-de time (the later of the two).
$overviewHandoffDate = "2016-03-24 16:31:24"

# Rows 7, 9-16 correspond to the files that received a new handoff in this
# run (row 8 - "In Translation" - is untouched).
$affectedRows = @(7, 9, 10, 11, 12, 13, 14, 15, 16)

$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $affectedRows) {
    $wsOverview.Cells.Item($r, 4).Value = $overviewHandoffDate
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $affectedRows) {
    $wsZhCn.Cells.Item($r, 5).Value = $zhcnHandoffTime
}

$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $affectedRows) {
    $wsDeDe.Cells.Item($r, 5).Value = $dedeHandoffTime
}
